$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 259, shifting existing
# rows 259:346 down to 260:347 (the last row becomes the new row 347).
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new weekly record.
$ws.Range("A259").Value = 8
$ws.Range("B259").Value = "Terminal La Palmera de La Serena"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 44588
$ws.Range("E259").Value = 4
$ws.Range("F259").Value = 100114001
$ws.Range("G259").Value = "Papa"
$ws.Range("H259").Value = "Asterix"
$ws.Range("I259").Value = "1a (cosecha)"
$ws.Range("J259").Value = 2600
$ws.Range("K259").Value = 9500
$ws.Range("L259").Value = 10000
$ws.Range("M259").Value = 9750
$ws.Range("N259").Value = "$/saco 25 kilos"
$ws.Range("O259").Value = "Provincia de Melipilla"
$ws.Range("P259").Value = 390
$ws.Range("Q259").Value = 25
$ws.Range("R259").Value = "Hortaliza"
